# Corrected excel sheets for application fix issues
#
# 1. NewLoanInput: update the test-case id string in B2, leave selection on B8
# 2. Edit Repayment Schedule: insert a new "waittopageload1" step (row 6,
#    formatted like the existing "waittopageload" row) before the Submit step,
#    and leave that new row selected/active.

$wb = $excel.ActiveWorkbook

# --- NewLoanInput sheet -----------------------------------------------
$ws1 = $wb.Worksheets.Item("NewLoanInput")
$ws1.Range("B2").Value = "2425-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-1-DATE-VAR-INST"
$ws1.Range("B8").Select() | Out-Null

# --- Edit Repayment Schedule sheet -------------------------------------
$ws5 = $wb.Worksheets.Item("Edit Repayment Schedule")
$ws5.Activate()

# Insert a new blank row above the old row 6 ("Submit" step), shifting
# everything below it down by one.
$ws5.Rows.Item(6).Insert() | Out-Null

# Populate the new row with the waittopageload1 step (name, wait-ms pair).
$ws5.Range("A6").Value = "waittopageload1"
$ws5.Range("B6").Value = 2000

# Match the formatting used by the existing "waittopageload" row (A3:B3).
$ws5.Range("B3").Copy() | Out-Null
$ws5.Range("B6").PasteSpecial(-4122) | Out-Null

$ws5.Range("A6:B6").Select() | Out-Null
